$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 (RM 232) - data completely removed, rows below shift up
$ws.Rows.Item(26).Delete()

# After the shift, the former "SC 92" row (was row 28) is now at row 27 - delete it too
$ws.Rows.Item(27).Delete()

# Update individual F-column values per the missing-data swap pattern
$ws.Range("F3").Value = 17.64
$ws.Range("F5").ClearContents()
$ws.Range("F21").Value = 16.58
$ws.Range("F23").ClearContents()
$ws.Range("F32").Value = 17.39
